$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Price (column D) values were refreshed by the scraper run.
# Cells are stored as text, so force text format before writing
# the new numeric-looking string and then restore the default
# (unstyled) cell style so no stray formatting is introduced.
$updates = [ordered]@{
    "D2" = "243.08"
    "D3" = "23.16"
    "D4" = "5.646"
    "D5" = "0.05825"
    "D6" = "3.408"
    "D7" = "6.487"
    "D9" = "0.7985"
    "D10" = "0.1468"
    "D11" = "0.07638"
    "D12" = "0.03249"
    "D13" = "0.03013"
    "D14" = "0.09232"
    "D15" = "0.001672"
    "D16" = "3.414"
    "D17" = "0.04764"
    "D18" = "0.0005995"
    "D20" = "0.001070"
    "D21" = "0.003829"
    "D22" = "0.0001501"
    "D23" = "3.693"
    "D24" = "2.210"
    "D25" = "0.3338"
    "D26" = "0.1254"
    "D27" = "0.0004003"
    "D40" = "0.04309"
    "D41" = "0.007017"
    "D42" = "0.1051"
    "D43" = "0.003386"
    "D44" = "0.008708"
    "D46" = "0.00005755"
    "D48" = "0.7859"
    "D49" = "0.1049"
}

foreach ($cell in $updates.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$cell]
    $rng.Style = "Normal"
}
